# Trade #19 closed at 2026-02-17 15:18:33 - unknown UNKNOWN +0.000%
#
# This script applies the same edit shown in the OOXML diff:
#  - Summary sheet: roll up totals to reflect the newly-closed trade
#  - Strategy Status sheet: update the MarketMaking row with new aggregates
#  - All Trades / MarketMaking sheets: append the new trade as row 20

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.8    # Current Capital
$summary.Range("B4").Value = -0.21     # Total P&L $
$summary.Range("B5").Value = -0.22     # Total P&L %
$summary.Range("B6").Value = 19        # Total Trades
$summary.Range("B8").Value = 9         # Losing Trades
$summary.Range("B9").Value = 26.32     # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.8       # Capital
$status.Range("D4").Value = 19         # Trades
$status.Range("E4").Value = -0.21      # P&L $
$status.Range("F4").Value = -0.2       # P&L %
$status.Range("G4").Value = 26.32      # Win Rate %

# ---------------------------------------------------------------
# 3. Append new trade row (#19) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$newRow = @{
    A = 19
    B = "2026-02-17"
    C = "15:18:27"
    D = "MarketMaking"
    E = "UP"
    F = 0.89
    G = 0.87
    H = "CLOSED"
    I = -2.2472
    J = -0.02
    K = 99.8
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A20").Value = $newRow.A

    # Date/time-looking text ("2026-02-17", "15:18:27") gets auto-coerced to a
    # real date/time serial by the COM layer's smart-entry heuristics. Force
    # plain text storage (matching the existing B/C columns in this sheet),
    # then reset the cell style back to Normal so no stray number-format
    # style sticks around on the new cells.
    $ws.Range("B20").NumberFormat = "@"
    $ws.Range("B20").Value = $newRow.B
    $ws.Range("B20").Style = "Normal"

    $ws.Range("C20").NumberFormat = "@"
    $ws.Range("C20").Value = $newRow.C
    $ws.Range("C20").Style = "Normal"
    $ws.Range("D20").Value = $newRow.D
    $ws.Range("E20").Value = $newRow.E
    $ws.Range("F20").Value = $newRow.F
    $ws.Range("G20").Value = $newRow.G
    $ws.Range("H20").Value = $newRow.H
    $ws.Range("I20").Value = $newRow.I
    $ws.Range("J20").Value = $newRow.J
    $ws.Range("K20").Value = $newRow.K
    $ws.Range("L20").Value = $newRow.L
    $ws.Range("M20").Value = $newRow.M
    $ws.Range("N20").Value = $newRow.N
    $ws.Range("O20").Value = $newRow.O
    $ws.Range("P20").Value = $newRow.P
    $ws.Range("Q20").Value = $newRow.Q
}
